$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.174.63"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "2.439.40"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'567.51"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Value = "'145.05"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("D9").Value = "'0.110"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("D16").Value = "62.145.47"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "2.440.00"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("E18").Value = "  -2.95%  "
$ws.Range("D19").Value = "'10.68"
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("D20").Value = "'319.31"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "'2.15"
$ws.Range("E22").Value = "  -3.21%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'9.81"
$ws.Range("E24").Value = "  +6.58%  "
$ws.Range("D25").Value = "'64.79"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").Value = "'638.17"
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").Value = "0.0₃0948"
$ws.Range("E28").Value = "  -5.55%  "
$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("D31").Value = "'7.80"
$ws.Range("E31").Value = "  -4.36%  "
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("D33").Value = "'0.131"
$ws.Range("E33").Value = "  -4.81%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("D36").Value = "'4.61"
$ws.Range("E36").Value = "  -3.88%  "
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").Value = "'0.363"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").Value = "'18.36"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("E40").Value = "  -5.48%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").Value = "'151.51"
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("D46").Value = "'15.30"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("D48").Value = "'0.599"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "'19.95"
$ws.Range("E49").Value = "  -4.44%  "
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").Value = "'0.0897"
$ws.Range("E51").Value = "  -2.57%  "
